# Applies the "Black Hat - fuja dessa pratica" section, inserted as a run of
# new list paragraphs right after the "Link site: https://moz.com/" paragraph
# (the last bullet of the Moz.com item) and before the trailing blank list
# paragraph that already closed the document.
$d = $word.ActiveDocument

# Locate the paragraph ending with the "Link site: https://moz.com/" hyperlink,
# which is the anchor point after which all new content must be inserted.
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Link site:*moz.com*") {
        $anchorIndex = $i
    }
}
if ($null -eq $anchorIndex) { throw "Anchor paragraph (Link site: https://moz.com/) not found" }
"Anchor paragraph index: $anchorIndex"

$insertAfter = $d.Paragraphs.Item($anchorIndex).Range

# Paragraph 1/20
$xml0 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">Black </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Hat</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> - fuja dessa prática</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
[void]$newPara.Range.InsertXML($xml0)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 1).Range

# Paragraph 2/20
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Os buscadores “premiam” os melhores sites com as primeiras páginas da SERP a partir do bom desempenho e boas práticas de cada site, mas eles também penalizam práticas ruins da mesma forma.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 2)
[void]$newPara.Range.InsertXML($xml1)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 2).Range

# Paragraph 3/20
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Há alguns anos, muitos profissionais de marketing digital chegaram à conclusão erroneamente que uma boa prática de SEO era a de “enganar o algoritmo” com uma série de táticas que tentavam burlar a inteligência de máquina dos buscadores e assim garantir as primeiras posições para o seu site nas páginas de pesquisa. E aí você se pergunta, se isso é antigo e se já sabemos que é uma prática condenável, existem profissionais que continuam fazendo?</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 3)
[void]$newPara.Range.InsertXML($xml2)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 3).Range

# Paragraph 4/20
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">A resposta infelizmente é sim! Esse conjunto de ações condenáveis pelos buscadores ganhou o nome de Black </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Hat</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> SEO, inspirado nos filmes de faroeste em que os vilões sempre usam o chapéu preto.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 4)
[void]$newPara.Range.InsertXML($xml3)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 4).Range

# Paragraph 5/20
$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Com o passar do tempo e todas as atualizações de algoritmo, o Google começou a identificar essas más práticas e deixou muito bem listado nas suas diretrizes o que não era aconselhável fazer. Separei uma lista do que é mais condenável e que deve ser evitado a todo custo no planejamento do seu marketing digital, isso porque muitos “gurus do marketing” continuam vendendo essas ideias como boas! Segue:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 5)
[void]$newPara.Range.InsertXML($xml4)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 5).Range

# Paragraph 6/20
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Uso excessivo de palavras-chave (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>keyword</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>stuffing</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 6)
[void]$newPara.Range.InsertXML($xml5)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 6).Range

# Paragraph 7/20
$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Se você estudou SEO no passado, já pode ter visto algo sobre isso: usar várias vezes a mesma palavra-chave para que o Google dê preferência ao seu conteúdo na hora de ranquear um site. Mas na verdade é o total oposto! Usar várias vezes o mesmo termo no corpo do texto ou conteúdo do seu site, além de ficar redundante e ruim para a leitura do usuário, também é uma prática condenável pelo Google. A plataforma já tem inteligência suficiente para entender sinônimos, por exemplo, então não precisa apelar para isso!</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 7)
[void]$newPara.Range.InsertXML($xml6)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 7).Range

# Paragraph 8/20
$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Conteúdo oculto</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 8)
[void]$newPara.Range.InsertXML($xml7)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 8).Range

# Paragraph 9/20
$xml8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Também já foi considerada uma boa estratégia usar conteúdo oculto nas páginas dos sites. Essa técnica era uma forma de inserir mais palavras-chave avulsas, outros termos de pesquisa e links de outros sites pensando apenas no robô de rastreamento dos buscadores, sem mostrar nada ao usuário. Esse conteúdo era omitido do usuário de formas bem simples como colocar o texto da mesma cor do fundo do site, alterar o texto para o tamanho da fonte em zero ou reposicionar o texto para fora da página usando CSS.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 9)
[void]$newPara.Range.InsertXML($xml8)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 9).Range

# Paragraph 10/20
$xml9 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Conteúdo duplicado/Plágio!</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 10)
[void]$newPara.Range.InsertXML($xml9)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 10).Range

# Paragraph 11/20
$xml10 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Como já falamos e muito aqui, pro Google conteúdo bom é conteúdo original! Não adianta pegar artigos de outros sites e replicar no seu blog achando que a frequência das postagens vai melhorar o seu posicionamento no ranking! Essa regra não diz nada sobre penalização, mas o Google não vai mostrar vários resultados idênticos ao usuário. Ao invés disso, ele prioriza o conteúdo original e oculta o restante que aparecer na SERP.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 11)
[void]$newPara.Range.InsertXML($xml10)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 11).Range

# Paragraph 12/20
$xml11 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Blog Spam</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 12)
[void]$newPara.Range.InsertXML($xml11)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 12).Range

# Paragraph 13/20
$xml12 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Anos atrás também era muito comum que os comentários de um fórum ou de um post de um blog fosse cheio de links para divulgar outras páginas ou até mesmo outros sites. Isso se feito repetidamente pode ser lido como spam e você pode receber penalizações do Google.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 13)
[void]$newPara.Range.InsertXML($xml12)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 13).Range

# Paragraph 14/20
$xml13 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Unrelated</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Keywords</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 14)
[void]$newPara.Range.InsertXML($xml13)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 14).Range

# Paragraph 15/20
$xml14 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Essa técnica se dá quando as palavras-chave utilizadas na otimização são de grande relevância, porém o conteúdo da página não é condizente. Essas palavras-chave vão desde nomes de celebridades até memes/bordões da Internet no meio de um conteúdo que não faz conexão com a página. Identificar isso é muito simples, pois não precisa ler muito do conteúdo para ver que as palavras são simplesmente encaixadas no texto sem o menor sentido. Muitas vezes elas são encontradas nos rodapés.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 15)
[void]$newPara.Range.InsertXML($xml14)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 15).Range

# Paragraph 16/20
$xml15 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>As punições do Google variam de acordo com a gravidade da infração cometida, o que pode ter relação com o número de vezes do mesmo erro, ou pelo peso de influência no desempenho do site. As possíveis punições são:</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 16)
[void]$newPara.Range.InsertXML($xml15)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 16).Range

# Paragraph 17/20
$xml16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Queda de 30 posições no ranking;</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 17)
[void]$newPara.Range.InsertXML($xml16)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 17).Range

# Paragraph 18/20
$xml17 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Queda de 50 posições;</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 18)
[void]$newPara.Range.InsertXML($xml17)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 18).Range

# Paragraph 19/20
$xml18 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Queda de 950 posições;</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 19)
[void]$newPara.Range.InsertXML($xml18)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 19).Range

# Paragraph 20/20
$xml19 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Banimento do Google.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertAfter.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 20)
[void]$newPara.Range.InsertXML($xml19)
$insertAfter = $d.Paragraphs.Item($anchorIndex + 20).Range

"Final paragraph count: $($d.Paragraphs.Count)"
